# New version read.xlsx and incorp NA fix v3.4 ALFAM2
$wb = $excel.ActiveWorkbook

# Rename the "Vars (hidden)" sheet to "Names"
$wsNames = $wb.Worksheets.Item("Vars (hidden)")
$wsNames.Name = "Names"

# On the Units sheet: relabel "Incorporation timing" -> "Time" and move the
# sheet's saved selection to C8
$wsUnits = $wb.Worksheets.Item("Units")
$wsUnits.Range("A8").Value = "Time"
$null = $wsUnits.Range("C8").Select()

# Make the renamed "Names" sheet the active tab (previously "Application" was active)
$null = $wsNames.Activate()
